# Two new price records were added to the "Papa" (Hortaliza, Femacal de
# La Calera) weekly consolidation sheet. They belong chronologically near
# the top of the existing block (rows 1221+), so insert two blank rows
# there (pushing the existing 1221-1282 data down to 1223-1284) and fill
# the two new rows with the new records.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 1221, shifting existing data down.
$ws.Rows.Item(1221).Insert()
$ws.Rows.Item(1221).Insert()

# --- New row 1221 ---------------------------------------------------
$ws.Range("A1221").Value = 3
$ws.Range("B1221").Value = "Femacal de La Calera"
$ws.Range("C1221").Value = "Coquimbo"
$ws.Range("D1221").Value = 45267
$ws.Range("E1221").Value = 5
$ws.Range("F1221").Value = 100114001
$ws.Range("G1221").Value = "Papa"
$ws.Range("H1221").Value = "Asterix"
$ws.Range("I1221").Value = "1a nueva(o)"
$ws.Range("J1221").Value = 160
$ws.Range("K1221").Value = 21000
$ws.Range("L1221").Value = 21000
$ws.Range("M1221").Value = 21000
$ws.Range("N1221").Value = "`$/saco 25 kilos"
$ws.Range("O1221").Value = "Provincia de Quillota"
$ws.Range("P1221").Value = 840
$ws.Range("Q1221").Value = 25
$ws.Range("R1221").Value = "Hortaliza"

# --- New row 1222 ---------------------------------------------------
$ws.Range("A1222").Value = 3
$ws.Range("B1222").Value = "Femacal de La Calera"
$ws.Range("C1222").Value = "Coquimbo"
$ws.Range("D1222").Value = 45267
$ws.Range("E1222").Value = 5
$ws.Range("F1222").Value = 100114001
$ws.Range("G1222").Value = "Papa"
$ws.Range("H1222").Value = "Rosara"
$ws.Range("I1222").Value = "1a nueva(o)"
$ws.Range("J1222").Value = 220
$ws.Range("K1222").Value = 20000
$ws.Range("L1222").Value = 21000
$ws.Range("M1222").Value = 20273
$ws.Range("N1222").Value = "`$/saco 25 kilos"
$ws.Range("O1222").Value = "Región de O'Higgins"
$ws.Range("P1222").Value = 811
$ws.Range("Q1222").Value = 25
$ws.Range("R1222").Value = "Hortaliza"

# Apply the same date number format ("s=2" style) used by the rest of
# the D column to the two new date cells, matching the column's style.
$ws.Range("D1221:D1222").NumberFormat = $ws.Range("D1223").NumberFormat
